$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 830.7714
$ws.Range("I33").Value = 559.375
$ws.Range("K33").Value = 559.375
$ws.Range("M33").Value = -330.375
$ws.Range("H98").Value = 903.23334
$ws.Range("I98").Value = 1033.174
$ws.Range("J98").Value = 476.2857
$ws.Range("K98").Value = 1033.174
$ws.Range("L98").Value = 476.2857
$ws.Range("M98").Value = 464.826
$ws.Range("N98").Value = -3472.2857
$ws.Range("H101").Value = 407.30768
$ws.Range("I101").Value = 358.63635
$ws.Range("J101").Value = 675
$ws.Range("K101").Value = 1075.90905
$ws.Range("L101").Value = 2025
$ws.Range("M101").Value = 546.09095
$ws.Range("N101").Value = -5269
$ws.Range("H122").Value = 903.23334
$ws.Range("I122").Value = 1033.174
$ws.Range("J122").Value = 476.2857
$ws.Range("K122").Value = 3099.522
$ws.Range("L122").Value = 1428.8571
$ws.Range("M122").Value = -649.5219999999999
$ws.Range("N122").Value = -6328.8571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 772.1
$ws.Range("I2").Value = 701.2857
$ws.Range("K2").Value = 701.2857
$ws.Range("M2").Value = -588.2857
$ws.Range("H61").Value = 3110
$ws.Range("I61").Value = 1911.0344
$ws.Range("K61").Value = 1911.0344
$ws.Range("M61").Value = -1699.0344
$ws.Range("H74").Value = 3680.037
$ws.Range("I74").Value = 1791.091
$ws.Range("J74").Value = 4978.6875
$ws.Range("K74").Value = 1791.091
$ws.Range("L74").Value = 4978.6875
$ws.Range("M74").Value = -917.0909999999999
$ws.Range("N74").Value = -6726.6875
$ws.Range("H77").Value = 3680.037
$ws.Range("I77").Value = 1791.091
$ws.Range("J77").Value = 4978.6875
$ws.Range("K77").Value = 8955.455
$ws.Range("L77").Value = 24893.4375
$ws.Range("M77").Value = -4587.455
$ws.Range("N77").Value = -33629.4375
$ws.Range("H110").Value = 1583.875
$ws.Range("I110").Value = 1457.762
$ws.Range("K110").Value = 1457.762
$ws.Range("M110").Value = 587.2380000000001
$ws.Range("H116").Value = 772.1
$ws.Range("I116").Value = 701.2857
$ws.Range("K116").Value = 701.2857
$ws.Range("M116").Value = 1592.7143
$ws.Range("H136").Value = 3110
$ws.Range("I136").Value = 1911.0344
$ws.Range("K136").Value = 5733.1032
$ws.Range("M136").Value = -3183.1032

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 772.1
$ws.Range("I3").Value = 701.2857
$ws.Range("K3").Value = 701.2857
$ws.Range("M3").Value = -587.2857
$ws.Range("H99").Value = 6596797
$ws.Range("I99").Value = 2266378.5
$ws.Range("J99").Value = 25001076
$ws.Range("K99").Value = 2266378.5
$ws.Range("L99").Value = 25001076
$ws.Range("M99").Value = -2264880.5
$ws.Range("N99").Value = -25004072
$ws.Range("H134").Value = 22990.244
$ws.Range("I134").Value = 29746.344
$ws.Range("J134").Value = 6100
$ws.Range("K134").Value = 89239.03200000001
$ws.Range("L134").Value = 18300
$ws.Range("M134").Value = -86704.03200000001
$ws.Range("N134").Value = -23370

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2911.875
$ws.Range("I31").Value = 2238.3235
$ws.Range("J31").Value = 3514.5264
$ws.Range("K31").Value = 2238.3235
$ws.Range("L31").Value = 3514.5264
$ws.Range("M31").Value = -1943.3235
$ws.Range("N31").Value = -4104.526400000001
$ws.Range("H34").Value = 2911.875
$ws.Range("I34").Value = 2238.3235
$ws.Range("J34").Value = 3514.5264
$ws.Range("K34").Value = 2238.3235
$ws.Range("L34").Value = 3514.5264
$ws.Range("M34").Value = -2036.3235
$ws.Range("N34").Value = -3918.5264
$ws.Range("H68").Value = 24939.166
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 24939.166
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 24939.166
$ws.Range("M68").ClearContents()
$ws.Range("N68").Value = -26437.166
$ws.Range("H71").Value = 24939.166
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 24939.166
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 74817.49800000001
$ws.Range("M71").ClearContents()
$ws.Range("N71").Value = -82305.49800000001
$ws.Range("H105").Value = 612.3036
$ws.Range("I105").Value = 604.1458
$ws.Range("J105").Value = 661.25
$ws.Range("K105").Value = 604.1458
$ws.Range("L105").Value = 661.25
$ws.Range("M105").Value = 1142.8542
$ws.Range("N105").Value = -4155.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 6980
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 6980
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 20940
$ws.Range("M76").ClearContents()
$ws.Range("N76").Value = -21706
$ws.Range("H79").Value = 6980
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 6980
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 20940
$ws.Range("M79").ClearContents()
$ws.Range("N79").Value = -23592

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 8781.833000000001
$ws.Range("I107").Value = 11375.111
$ws.Range("K107").Value = 11375.111
$ws.Range("M107").Value = -9455.111000000001
$ws.Range("H132").Value = 3968
$ws.Range("I132").Value = 3944.3845
$ws.Range("K132").Value = 11833.1535
$ws.Range("M132").Value = -9303.1535

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2734.4443
$ws.Range("I7").Value = 3043.75
$ws.Range("J7").Value = 2487
$ws.Range("K7").Value = 3043.75
$ws.Range("L7").Value = 2487
$ws.Range("M7").Value = -2931.75
$ws.Range("N7").Value = -2711
$ws.Range("H16").Value = 3745.0952
$ws.Range("I16").Value = 3263.7222
$ws.Range("J16").Value = 6633.3335
$ws.Range("K16").Value = 3263.7222
$ws.Range("L16").Value = 6633.3335
$ws.Range("M16").Value = -3093.7222
$ws.Range("N16").Value = -6973.3335
$ws.Range("H64").Value = 23975
$ws.Range("J64").Value = 23975
$ws.Range("L64").Value = 23975
$ws.Range("N64").Value = -24425
$ws.Range("H67").Value = 23975
$ws.Range("J67").Value = 23975
$ws.Range("L67").Value = 23975
$ws.Range("N67").Value = -25535
$ws.Range("H68").Value = 43480180
$ws.Range("I68").Value = 62501350
$ws.Range("J68").Value = 3214.2856
$ws.Range("K68").Value = 62501350
$ws.Range("L68").Value = 3214.2856
$ws.Range("M68").Value = -62500601
$ws.Range("N68").Value = -4712.2856
$ws.Range("H71").Value = 43480180
$ws.Range("I71").Value = 62501350
$ws.Range("J71").Value = 3214.2856
$ws.Range("K71").Value = 312506750
$ws.Range("L71").Value = 16071.428
$ws.Range("M71").Value = -312503006
$ws.Range("N71").Value = -23559.428
$ws.Range("H126").Value = 2734.4443
$ws.Range("I126").Value = 3043.75
$ws.Range("J126").Value = 2487
$ws.Range("K126").Value = 9131.25
$ws.Range("L126").Value = 7461
$ws.Range("M126").Value = -6661.25
$ws.Range("N126").Value = -12401
$ws.Range("H132").Value = 6883
$ws.Range("I132").Value = 2021
$ws.Range("J132").Value = 16607
$ws.Range("K132").Value = 6063
$ws.Range("L132").Value = 49821
$ws.Range("M132").Value = -3533
$ws.Range("N132").Value = -54881
$ws.Range("H136").Value = 4608.1113
$ws.Range("I136").Value = 2538.44
$ws.Range("J136").Value = 7195.2
$ws.Range("K136").Value = 7615.32
$ws.Range("L136").Value = 21585.6
$ws.Range("M136").Value = -5065.32
$ws.Range("N136").Value = -26685.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H122").Value = 38382.89
$ws.Range("I122").Value = 39666.848
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 119000.544
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -116550.544
$ws.Range("N122").Value = -19900
$ws.Range("H136").Value = 2282.0667
$ws.Range("I136").Value = 1813.4762
$ws.Range("J136").Value = 2692.0833
$ws.Range("K136").Value = 5440.4286
$ws.Range("L136").Value = 8076.249899999999
$ws.Range("M136").Value = -2890.4286
$ws.Range("N136").Value = -13176.2499
